$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column values that look numeric stay as text, matching the source data.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.953.72"
$ws.Range("E2").Value = "  +0.39%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.576.94"
$ws.Range("E3").Value = "  +2.16%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.32"
$ws.Range("E5").Value = "  +0.41%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.87"
$ws.Range("E6").Value = "  +4.15%  "

$ws.Range("E7").Value = "  -0.15%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("E9").Value = "  +0.57%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.38"
$ws.Range("E10").Value = "  +0.71%  "

$ws.Range("E11").Value = "  +0.40%  "

$ws.Range("E12").Value = "  -0.27%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.972.97"
$ws.Range("E13").Value = "  +2.21%  "

$ws.Range("E14").Value = "  -0.49%  "

$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.75"
$ws.Range("E15").Value = "  +3.12%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.576.22"
$ws.Range("E16").Value = "  +4.20%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.845"
$ws.Range("E17").Value = "  -1.03%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.026.49"
$ws.Range("E18").Value = "  +0.42%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.85"
$ws.Range("E19").Value = "  +1.49%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.68"
$ws.Range("E20").Value = "  -1.50%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0969"
$ws.Range("E21").Value = "  +0.58%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.55"
$ws.Range("E22").Value = "  -0.27%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "250.03"
$ws.Range("E23").Value = "  -1.46%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.98"
$ws.Range("E24").Value = "  +1.05%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.08"
$ws.Range("E25").Value = "  -0.01%  "

$ws.Range("E26").Value = "  +1.81%  "

$ws.Range("E27").Value = "  +0.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.39"
$ws.Range("E28").Value = "  -1.14%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.61"
$ws.Range("E29").Value = "  -0.77%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.32"
$ws.Range("E30").Value = "  -0.51%  "

$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "157.90"
$ws.Range("E31").Value = "  +0.23%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.84"
$ws.Range("E32").Value = "  -1.92%  "

$ws.Range("E33").Value = "  +3.82%  "

$ws.Range("E34").Value = "  -1.18%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0804"
$ws.Range("E35").Value = "  +3.04%  "

$ws.Range("E36").Value = "  -0.36%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.80"
$ws.Range("E37").Value = "  -3.72%  "

$ws.Range("E38").Value = "  +10.08%  "

$ws.Range("E39").Value = "  +1.34%  "

$ws.Range("E40").Value = "  +0.22%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "23.49"
$ws.Range("E41").Value = "  +1.10%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.14"
$ws.Range("E42").Value = "  +9.14%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0305"
$ws.Range("E43").Value = "  -0.28%  "

$ws.Range("E44").Value = "  -0.07%  "

$ws.Range("E45").Value = "  -2.15%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.002.94"
$ws.Range("E46").Value = "  -2.29%  "

$ws.Range("E47").Value = "  +0.20%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.822.30"
$ws.Range("E48").Value = "  +2.07%  "

$ws.Range("E49").Value = "  +2.78%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "75.03"
$ws.Range("E50").Value = "  -0.95%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "81.87"
$ws.Range("E51").Value = "  -4.22%  "
